{"js": "// Update the instructor's Tuesday office-hours time, and fill in the\n// previously-TBD TA office hours time/location.\n//\n// \"Tuesday 2:15 - 3:30pm\"  ->  \"Tuesday 3:30 - 4:45pm\"\n// \"TBD\" (TA Office Hours)  ->  \"Mondays, Wednesdays 7-9pm in WNS 216\"\n\nconst body = context.document.body;\n\n// 1) Instructor office hours: Tuesday slot moves from 2:15-3:30pm to 3:30-4:45pm.\nconst tuesdayResults = body.search(\"Tuesday 2:15 - 3:30pm\", { matchCase: true, matchWholeWord: false });\ntuesdayResults.load(\"items\");\nawait context.sync();\n\nif (tuesdayResults.items.length > 0) {\n  // insertText(..., replace) swaps the matched range's text in place while\n  // keeping the run's existing formatting (font, size, shading, etc.).\n  tuesdayResults.items[0].insertText(\"Tuesday 3:30 - 4:45pm\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) TA office hours: replace the \"TBD\" placeholder with the actual schedule.\nconst tbdResults = body.search(\"TBD\", { matchCase: true, matchWholeWord: true });\ntbdResults.load(\"items\");\nawait context.sync();\n\nif (tbdResults.items.length > 0) {\n  tbdResults.items[0].insertText(\"Mondays, Wednesdays 7-9pm in WNS 216\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the instructor's Tuesday office-hours time, and fill in the\n# previously-TBD TA office hours time/location.\n#\n# \"Tuesday 2:15 - 3:30pm\"  ->  \"Tuesday 3:30 - 4:45pm\"\n# \"TBD\" (TA Office Hours)  ->  \"Mondays, Wednesdays 7-9pm in WNS 216\"\n\n$d = $word.ActiveDocument\n\n# 1) Instructor office hours: Tuesday slot moves from 2:15-3:30pm to 3:30-4:45pm.\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\n    \"Tuesday 2:15 - 3:30pm\", # FindText\n    $false,                  # MatchCase\n    $false,                  # MatchWholeWord\n    $false,                  # MatchWildcards\n    $false,                  # MatchSoundsLike\n    $false,                  # MatchAllWordForms\n    $true,                   # Forward\n    1,                       # Wrap (wdFindContinue)\n    $false,                  # Format\n    \"Tuesday 3:30 - 4:45pm\", # ReplaceWith\n    2                        # Replace (wdReplaceAll)\n) | Out-Null\n\n# 2) TA office hours: replace the \"TBD\" placeholder with the actual schedule.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\n    \"TBD\",                                    # FindText\n    $false,                                   # MatchCase\n    $true,                                    # MatchWholeWord\n    $false,                                   # MatchWildcards\n    $false,                                   # MatchSoundsLike\n    $false,                                   # MatchAllWordForms\n    $true,                                    # Forward\n    1,                                        # Wrap (wdFindContinue)\n    $false,                                   # Format\n    \"Mondays, Wednesdays 7-9pm in WNS 216\",   # ReplaceWith\n    2                                         # Replace (wdReplaceAll)\n) | Out-Null\n"}
